# "Colocando header nos graficos"
# Add a header label in column A (row 1) for each data sheet, fix a few
# accented Portuguese labels, strip the bold/bordered header style from the
# row-label cells (A2:A12, etc.) since they are no longer header cells, and
# correct a couple of values/labels on the last two sheets. Also remove the
# now-unused "Teto" row on the emissions sheet.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122
$xlPasteValues = -4163

# ---------------------------------------------------------------------
# Sheets 1-4 share an identical layout/fix-up:
#   Potencia Acumulada - SIN (MW)
#   Geracao Periodo Medio (MWMed)
#   Atendimento a Ponta(MW)
#   Potencia Incremental - SIN(MW)
# ---------------------------------------------------------------------
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # New header cell for the technology/source column, styled like the
    # other header cells on row 1 (bold, centered, bordered).
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial($xlPasteFormats)
    $ws.Application.CutCopyMode = $false

    # Fix accented labels.
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."

    # The row-label cells are no longer header cells, so drop the bold
    # bordered header style back to the default "Normal" style.
    $ws.Range("A2:A12").Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet 5: Emissoes Totais (MtCO2eq)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("A1").Value = "Período"
$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial($xlPasteFormats)
$ws5.Application.CutCopyMode = $false

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A2:A3").Style = "Normal"

# Remove the now-unused "Teto" row entirely.
$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6: Custo Total (bilhões de R$)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial($xlPasteFormats)
$ws6.Application.CutCopyMode = $false

# Re-label the value column header as "2015" (kept as text, like the
# other sheets' year headers) while preserving the existing header style.
# Copying the text value from sheet 1's "2015" header cell avoids Excel
# auto-converting a literal "2015" string assignment into a number.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial($xlPasteValues)
$ws6.Application.CutCopyMode = $false

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 568
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
$ws6.Range("A2:A3").Style = "Normal"
